# Generate Report for Handoff
# Updates the "Latest HO Xliff Generate Date" / "Latest Handoff Datetime"
# timestamps and the "Priority" column for the rows that were just
# (re-)handed off, on the Overview, zh-cn and de-de sheets.

$wb = $excel.ActiveWorkbook

$rows = @(7, 9, 10, 11, 12, 13)

# --- Overview sheet: column G = "Latest HO Xliff Generate Date" ---
$overview = $wb.Worksheets.Item("Overview")
foreach ($r in $rows) {
    $overview.Range("G$r").Value = "2016-09-05 00:25:35"
}

# --- zh-cn sheet: column H = "Latest Handoff Datetime", column E = "Priority" ---
$zhcn = $wb.Worksheets.Item("zh-cn")
foreach ($r in $rows) {
    $zhcn.Range("H$r").Value = "2016-09-05 00:25:30"
    $zhcn.Range("E$r").Value = "ht"
}

# --- de-de sheet: column H = "Latest Handoff Datetime", column E = "Priority" ---
$dede = $wb.Worksheets.Item("de-de")
foreach ($r in $rows) {
    $dede.Range("H$r").Value = "2016-09-05 00:25:35"
    $dede.Range("E$r").Value = "ht"
}
